$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal: after the existing last paragraph ("Role 2:"), append:
#   - a bold "TAGS:" header paragraph (same style as the document's
#     other bold "Label:" headers)
#   - nine plain (non-bold) paragraphs, one per tag line
# ------------------------------------------------------------------

# 1) Append the bold "TAGS:" header right after "Role 2:". Doing this
#    as the very next paragraph after an already-bold paragraph makes
#    Word carry the bold run/paragraph formatting forward cleanly
#    (matching how every other "Label:" header in this document is
#    represented).
$roleRange = $d.Paragraphs.Last.Range
$roleRange.InsertParagraphAfter()
$tagsPara = $d.Paragraphs.Last
$tagsPara.Range.Text = "TAGS:"

# 2) Build the nine tag-line paragraphs someplace that is NOT bold, so
#    they come out as plain paragraphs with no direct character
#    formatting (mirrors other plain paragraphs already in the file).
#    The paragraph right after the document's second paragraph
#    ("cpb-aacip-...") is a convenient, guaranteed non-bold anchor.
$anchorPara = $d.Paragraphs.Item(2)
$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphAfter()
$insertPoint = $d.Paragraphs.Item(3).Range
$insertPoint.Collapse(0)

$lines = @(
    "1 Cikuliurun, Tugeq -- Ice Spud, Ice Pick, Ice Chisel",
    "1 Ciku -- Ice",
    "1 Pellaalleq -- Getting Lost",
    "1 Quyayaraq -- Thanksgiving",
    "1 Kiuyaraq, Kiutaaryaraq -- Talking Back",
    "1 Qanruyutet, Qaneryarat -- Traditional Wisdom, Wise Words",
    "1 Qanpautevkenaki -- Do Not Shout at Them",
    "1 Elliraat -- Orphans",
    "1 Nerangnaqsaraq, Yuungnaqsaraq -- Subsistence"
)

for ($i = 0; $i -lt $lines.Count; $i++) {
    $p = $d.Paragraphs.Item(3 + $i)
    $p.Range.Text = $lines[$i]
    if ($i -lt $lines.Count - 1) {
        $p.Range.InsertParagraphAfter()
    }
}

# 3) Cut that freshly-built, cleanly-formatted block of nine
#    paragraphs out of its temporary location.
$blockStart = $d.Paragraphs.Item(3).Range.Start
$blockEndPara = $d.Paragraphs.Item(3 + $lines.Count - 1)
$blockEnd = $blockEndPara.Range.End
$blockRange = $d.Range($blockStart, $blockEnd)
$blockRange.Cut()

# 4) Paste it back in immediately after "TAGS:" (now once again the
#    last paragraph in the document), after opening up a fresh
#    paragraph for the pasted content to land in.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$pasteRange = $d.Content
$pasteRange.Collapse(0)
$pasteRange.Paste()

# 5) The paste leaves one extra empty trailing paragraph behind (the
#    one opened in step 4) -- remove it so the document ends exactly
#    on the last tag line, as it did on "Role 2:" before this edit.
$lastPara = $d.Paragraphs.Last
if ($lastPara.Range.Text.Trim() -eq "") {
    $lastPara.Range.Delete()
}
